# Master Data refresh (16th May) - add 3 new reg_center_user_machine rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows appended after the existing data (rows 2-33), continuing the
# usr_id/machine_id sequence (110033-110035) for regcntr_id/lang_code/
# is_active/cr_by/cr_dtimes matching the existing pattern.
$newRows = @(
    @(10005, 110033, 10005, "eng", $true, "superadmin", "now()"),
    @(10005, 110034, 10005, "eng", $true, "superadmin", "now()"),
    @(10005, 110035, 10005, "eng", $true, "superadmin", "now()")
)

$rowIndex = 34
foreach ($rowData in $newRows) {
    $ws.Cells.Item($rowIndex, 1).Value = $rowData[0]
    $ws.Cells.Item($rowIndex, 2).Value = $rowData[1]
    $ws.Cells.Item($rowIndex, 3).Value = $rowData[2]
    $ws.Cells.Item($rowIndex, 4).Value = $rowData[3]
    $ws.Cells.Item($rowIndex, 5).Value = $rowData[4]
    $ws.Cells.Item($rowIndex, 6).Value = $rowData[5]
    $ws.Cells.Item($rowIndex, 7).Value = $rowData[6]
    $rowIndex++
}

# Select the row below the last data row through the end of the sheet,
# matching the selection left behind after entering the new rows.
$ws.Range($ws.Rows.Item(37), $ws.Rows.Item(1048576)).Select() | Out-Null
